$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.944.44"
$ws.Range("E2").Value = "  -0.74%  "
$ws.Range("D3").Value = "2.499.20"
$ws.Range("E3").Value = "  +0.93%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'541.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.65%  "
$ws.Range("D6").Value = "'143.69"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.89%  "
$ws.Range("D7").Value = "'0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  +0.45%  "
$ws.Range("D9").Value = "2.521.48"
$ws.Range("E9").Value = "  +1.88%  "
$ws.Range("E10").Value = "  +1.07%  "
$ws.Range("E11").Value = "  +0.51%  "
$ws.Range("E12").Value = "  +4.75%  "
$ws.Range("D13").Value = "'0.355"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.90%  "
$ws.Range("D14").Value = "2.945.89"
$ws.Range("E14").Value = "  +1.70%  "
$ws.Range("D15").Value = "'23.50"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.76%  "
$ws.Range("D16").Value = "58.871.71"
$ws.Range("E16").Value = "  -0.80%  "
$ws.Range("E17").Value = "  +1.14%  "
$ws.Range("D18").Value = "2.518.19"
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("D19").Value = "'11.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.32%  "
$ws.Range("E20").Value = "  -1.49%  "
$ws.Range("D21").Value = "'324.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("D22").Value = "'0.997"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.06%  "
$ws.Range("E23").Value = "  +0.88%  "
$ws.Range("D24").Value = "'62.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E25").Value = "  -5.07%  "
$ws.Range("E26").Value = "  +0.84%  "
$ws.Range("D27").Value = "2.624.66"
$ws.Range("E27").Value = "  +1.78%  "
$ws.Range("E28").Value = "  +1.72%  "
$ws.Range("D29").Value = "'7.83"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.13%  "
$ws.Range("D30").Value = "0.0₃0775"
$ws.Range("E30").Value = "  +0.28%  "
$ws.Range("E31").Value = "  -0.51%  "
$ws.Range("E32").Value = "  -2.56%  "
$ws.Range("E33").Value = "  -5.77%  "
$ws.Range("D34").Value = "'0.997"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'1.45"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.38%  "
$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D36").Value = "'156.41"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.52%  "
$ws.Range("D37").Value = "'18.63"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.35%  "
$ws.Range("E38").Value = "  -5.55%  "
$ws.Range("E39").Value = "  -9.53%  "
$ws.Range("D40").Value = "'5.67"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.96%  "
$ws.Range("D41").Value = "'36.87"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.27%  "
$ws.Range("D42").Value = "'296.12"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.27%  "
$ws.Range("D43").Value = "'3.69"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.65%  "
$ws.Range("E44").Value = "  -2.94%  "
$ws.Range("D45").Value = "'0.995"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.30%  "
$ws.Range("D46").Value = "'0.597"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.85%  "
$ws.Range("D47").Value = "'10.77"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.33%  "
$ws.Range("E48").Value = "  -1.17%  "
$ws.Range("D49").Value = "'122.86"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.63%  "
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("E51").Value = "  -0.38%  "
